$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New input rows describing waiting time ("espera") per trip distance.
# New shared strings must be created in the order: "km", "espera (€)",
# "costo hora de espera/km" so they land at shared-string indices 7,8,9
# respectively (matching the target workbook).
# ---------------------------------------------------------------------

# Row 18: distance header (km), bold
$ws.Range("A18").Value = "km"

$kms = @(10,11,12,13,14,15,16,17,18,19,20,21)
$col = 2
foreach ($v in $kms) {
    $ws.Cells.Item(18, $col).Value = $v
    $col = $col + 1
}
$ws.Range("A18:M18").Font.Bold = $true

# Row 19: waiting cost per distance bracket (€)
$ws.Range("A19").Value = "espera (€)"

$esperas = @(15,14,14,13,13,12,12,11,11,10,10,10)
$fillColor = 3506772
$col = 2
foreach ($v in $esperas) {
    $cell = $ws.Cells.Item(19, $col)
    $cell.Value = $v
    $cell.Interior.Color = $fillColor
    $col = $col + 1
}

# Row 15: new coefficient row, alongside the existing "coeficiente HP" /
# "coeficiente Asientos" rows (13 and 14).
$ws.Range("A15").Value = "costo hora de espera/km"
$ws.Range("B15").Interior.Color = $fillColor

# Column A is very slightly wider in the updated workbook.
$ws.Columns.Item(1).ColumnWidth = 23.67

# Selection moved to P20 in the saved file.
$null = $ws.Range("P20").Select()
